$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "SW_MATM SW_DAI"
$ws.Range("C3").Value = "SSH SEC_LOGIN"
$ws.Range("C5").Value = "SW_DAI"
$ws.Range("C10").Value = "DTP"
$ws.Range("C11").Value = "MACNOTIFY PORT_SECURITY"

$ws.Range("C13").Select()
